$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.318.53"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.032.40"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.74"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.01"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.029.72"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("E11").Value = "  -1.72%  "

$ws.Range("E12").Value = "  +5.42%  "

$ws.Range("E13").Value = "  -2.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.47"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +4.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.268.25"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.534.08"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("E18").Value = "  +4.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.59"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +19.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.031.85"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "475.84"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.21"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("E25").Value = "  +4.04%  "

$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -4.20%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.44"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0998"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  -5.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +5.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.27"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  +3.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.994"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.47"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +9.83%  "

$ws.Range("E39").Value = "  -5.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.56"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -0.28%  "

$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("E43").Value = "  -5.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.63"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0360"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "380.95"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  -5.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.721.13"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.84"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("E51").Value = "  +3.63%  "
